# Scheduled-runner sync: refresh Universalis market-price snapshots
# (currentAveragePrice / NQ / HQ + derived leve profit columns) across
# every crafter sheet in the Behemoth Profits workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 33: "Glazed and Confused" (Clear Glass Lens)
$ws.Range("H33").Value = 977.4
$ws.Range("I33").Value = 866.65216
$ws.Range("J33").Value = 2251
$ws.Range("K33").Value = 866.65216
$ws.Range("L33").Value = 2251
$ws.Range("M33").Value = -637.65216
$ws.Range("N33").Value = -2709

# Row 62: "The Mustache Suits Him" (Enchanted Mythrite Ink)
$ws.Range("H62").Value = 1858.2
$ws.Range("I62").Value = 1766.6666
$ws.Range("J62").Value = 1995.5
$ws.Range("K62").Value = 1766.6666
$ws.Range("L62").Value = 1995.5
$ws.Range("M62").Value = -1142.6666
$ws.Range("N62").Value = -3243.5

# Row 65: "Forgery of Convenience (L)" (Enchanted Mythrite Ink)
$ws.Range("H65").Value = 1858.2
$ws.Range("I65").Value = 1766.6666
$ws.Range("J65").Value = 1995.5
$ws.Range("K65").Value = 8833.333000000001
$ws.Range("L65").Value = 9977.5
$ws.Range("M65").Value = -5713.333000000001
$ws.Range("N65").Value = -16217.5

# Row 70: "Consecrating Congregation" (Holy Water)
$ws.Range("H70").Value = 7499.5
$ws.Range("I70").Value = 4999
$ws.Range("J70").Value = 10000
$ws.Range("K70").Value = 14997
$ws.Range("L70").Value = 30000
$ws.Range("M70").Value = -14727
$ws.Range("N70").Value = -30540

# Row 73: "Curbing the Contagion (L)" (Holy Water)
$ws.Range("H73").Value = 7499.5
$ws.Range("I73").Value = 4999
$ws.Range("J73").Value = 10000
$ws.Range("K73").Value = 14997
$ws.Range("L73").Value = 30000
$ws.Range("M73").Value = -14061
$ws.Range("N73").Value = -31872

# Row 76: "Warding Off Temptation" (Enchanted Hardsilver Ink)
$ws.Range("H76").Value = 4556.7144
$ws.Range("I76").Value = 4074.5
$ws.Range("J76").Value = 5199.6665
$ws.Range("K76").Value = 4074.5
$ws.Range("L76").Value = 5199.6665
$ws.Range("M76").Value = -3759.5
$ws.Range("N76").Value = -5829.6665

# Row 79: "The Garden of Arcane Delights (L)" (Enchanted Hardsilver Ink)
$ws.Range("H79").Value = 4556.7144
$ws.Range("I79").Value = 4074.5
$ws.Range("J79").Value = 5199.6665
$ws.Range("K79").Value = 4074.5
$ws.Range("L79").Value = 5199.6665
$ws.Range("M79").Value = -2982.5
$ws.Range("N79").Value = -7383.6665

# Row 103: "Let Loose the Juice" (Persimmon Tannin)
$ws.Range("H103").Value = 1000
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 1000
$ws.Range("K103").Value = 0
$ws.Range("L103").ClearContents()
$ws.Range("M103").Value = 3000
$ws.Range("N103").Value = -4172


$ws = $wb.Worksheets.Item("ARM")

# Row 39: "Aurochs Star" (Bull Hoplon)
$ws.Range("H39").Value = 4470
$ws.Range("I39").Value = 4333
$ws.Range("K39").Value = 4333
$ws.Range("M39").Value = -3813

# Row 74: "As the Bolt Flies" (Titanium Nugget)
$ws.Range("H74").Value = 12510178
$ws.Range("I74").Value = 27779360
$ws.Range("J74").Value = 17211.727
$ws.Range("K74").Value = 27779360
$ws.Range("L74").Value = 17211.727
$ws.Range("M74").Value = -27778486
$ws.Range("N74").Value = -18959.727

# Row 77: "Heavy Metal Banned (L)" (Titanium Nugget)
$ws.Range("H77").Value = 12510178
$ws.Range("I77").Value = 27779360
$ws.Range("J77").Value = 17211.727
$ws.Range("K77").Value = 138896800
$ws.Range("L77").Value = 86058.63499999999
$ws.Range("M77").Value = -138892432
$ws.Range("N77").Value = -94794.63499999999

# Row 122: "Haste for High Durium" (High Durium Nugget)
$ws.Range("H122").Value = 1999
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1999
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").Value = 5997
$ws.Range("N122").Value = -10897

# Row 132: "Don't Bore Me, Ore Me" (Mountain Chromite Ingot)
$ws.Range("H132").Value = 5305.448
$ws.Range("I132").Value = 2212.4634
$ws.Range("K132").Value = 6637.3902
$ws.Range("M132").Value = -4107.3902


$ws = $wb.Worksheets.Item("BSM")

# Row 20: "Smelt and Dealt" (Iron Ingot)
$ws.Range("H20").Value = 6128.6665
$ws.Range("I20").Value = 6128.6665
$ws.Range("K20").Value = 6128.6665
$ws.Range("M20").Value = -5881.6665

# Row 22: "Riveting Run" (Iron Rivets)
$ws.Range("H22").Value = 520.5
$ws.Range("I22").Value = 281
$ws.Range("J22").Value = 999.5
$ws.Range("K22").Value = 281
$ws.Range("L22").Value = 999.5
$ws.Range("M22").Value = -108
$ws.Range("N22").Value = -1345.5

# Row 25: "Tools of the Trade" (Iron Doming Hammer)
$ws.Range("H25").Value = 399.66666
$ws.Range("I25").Value = 399.66666
$ws.Range("K25").Value = 399.66666
$ws.Range("M25").Value = -164.66666

# Row 94: "High Steal" (High Steel Nugget)
$ws.Range("H94").Value = 1013.1923
$ws.Range("I94").Value = 973.72
$ws.Range("J94").Value = 2000
$ws.Range("K94").Value = 973.72
$ws.Range("L94").Value = 2000
$ws.Range("M94").Value = -522.72
$ws.Range("N94").Value = -2902

# Row 105: "Ingot to Wing It" (Molybdenum Ingot)
$ws.Range("H105").Value = 2838
$ws.Range("I105").Value = 2554.75
$ws.Range("J105").Value = 2999.8572
$ws.Range("K105").Value = 2554.75
$ws.Range("L105").Value = 2999.8572
$ws.Range("M105").Value = -807.75
$ws.Range("N105").Value = -6493.8572

# Row 134: "Ruthenium Supremium" (Ruthenium Ingot)
$ws.Range("H134").Value = 36926.465
$ws.Range("I134").Value = 1256.5555
$ws.Range("K134").Value = 3769.6665
$ws.Range("M134").Value = -1234.6665


$ws = $wb.Worksheets.Item("CRP")

# Row 35: "Storm of Swords" (Elm Macuahuitl)
$ws.Range("H35").Value = 362.5
$ws.Range("I35").Value = 362.5
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 362.5
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -68.5

# Row 56: "Trident and Error" (Cobalt Trident)
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()


$ws = $wb.Worksheets.Item("CUL")

# Row 6: "Meat-lover's Special" (Marmot Steak)
$ws.Range("H6").Value = 454.57144
$ws.Range("I6").Value = 530
$ws.Range("J6").Value = 2
$ws.Range("K6").Value = 1590
$ws.Range("L6").Value = 6
$ws.Range("M6").Value = -1477
$ws.Range("N6").Value = -232


$ws = $wb.Worksheets.Item("GSM")

# Row 11: "A Ringing Success" (Copper Ring)
$ws.Range("H11").Value = 9159090
$ws.Range("I11").Value = 6009875
$ws.Range("J11").Value = 17032126
$ws.Range("K11").Value = 6009875
$ws.Range("L11").Value = 17032126
$ws.Range("M11").Value = -6009736
$ws.Range("N11").Value = -17032404

# Row 62: "The Goggles, They Do Naught" (Mythrite Goggles of Gathering)
$ws.Range("H62").Value = 99957.664
$ws.Range("J62").Value = 99957.664
$ws.Range("L62").Value = 99957.664
$ws.Range("N62").Value = -101329.664

# Row 65: "Peril Never Wore Safety Goggles (L)" (Mythrite Goggles of Gathering)
$ws.Range("H65").Value = 99957.664
$ws.Range("J65").Value = 99957.664
$ws.Range("L65").Value = 299872.992
$ws.Range("N65").Value = -306736.992

# Row 70: "Sky Is the Limit" (Mythrite Ingot)
$ws.Range("H70").Value = 4101.1665
$ws.Range("I70").Value = 4101.1665
$ws.Range("K70").Value = 4101.1665
$ws.Range("M70").Value = -3831.1665

# Row 73: "Hulls of Broken Dreams (L)" (Mythrite Ingot)
$ws.Range("H73").Value = 4101.1665
$ws.Range("I73").Value = 4101.1665
$ws.Range("K73").Value = 4101.1665
$ws.Range("M73").Value = -3165.1665

# Row 122: "Awarding Academic Excellence" (Ametrine)
$ws.Range("H122").Value = 1999.625
$ws.Range("I122").Value = 1959.4
$ws.Range("K122").Value = 5878.200000000001
$ws.Range("M122").Value = -3428.200000000001

# Row 126: "Gold Rush Order" (Phrygian Gold Ingot)
$ws.Range("H126").Value = 1662
$ws.Range("I126").Value = 1729.091
$ws.Range("J126").Value = 1293
$ws.Range("K126").Value = 5187.272999999999
$ws.Range("L126").Value = 3879
$ws.Range("M126").Value = -2717.272999999999
$ws.Range("N126").Value = -8819

# Row 132: "On Board for Lar" (Lar Ingot)
$ws.Range("H132").Value = 71431416
$ws.Range("I132").Value = 83335920
$ws.Range("J132").Value = 4394.5
$ws.Range("K132").Value = 250007760
$ws.Range("L132").Value = 13183.5
$ws.Range("M132").Value = -250005230
$ws.Range("N132").Value = -18243.5


$ws = $wb.Worksheets.Item("LTW")

# Row 7: "Tan Before the Ban" (Leather)
$ws.Range("H7").Value = 57610.105
$ws.Range("I7").Value = 4172.75
$ws.Range("J7").Value = 149217
$ws.Range("K7").Value = 4172.75
$ws.Range("L7").Value = 149217
$ws.Range("M7").Value = -4060.75
$ws.Range("N7").Value = -149441

# Row 103: "Security Breeches" (Marid Leather Breeches of Maiming)
$ws.Range("H103").Value = 35522.25
$ws.Range("J103").Value = 35522.25
$ws.Range("L103").Value = 35522.25
$ws.Range("N103").Value = -37866.25

# Row 111: "Glove Me Tender" (Gliderskin Gloves of Striking)
$ws.Range("H111").Value = 134000
$ws.Range("J111").Value = 134000
$ws.Range("L111").Value = 134000
$ws.Range("N111").Value = -142180

# Row 122: "Hell on Leather" (Gaja Leather)
$ws.Range("H122").Value = 5412.2173
$ws.Range("I122").Value = 4398.875
$ws.Range("K122").Value = 13196.625
$ws.Range("M122").Value = -10746.625

# Row 126: "Battered Books" (Saiga Leather)
$ws.Range("H126").Value = 57610.105
$ws.Range("I126").Value = 4172.75
$ws.Range("J126").Value = 149217
$ws.Range("K126").Value = 12518.25
$ws.Range("L126").Value = 447651
$ws.Range("M126").Value = -10048.25
$ws.Range("N126").Value = -452591


$ws = $wb.Worksheets.Item("WVR")

# Row 12: "This Is Why You Can't Have Nice Things" (Amateur's Breeches)
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
